# Update "想去人数" (interest count) values across all four sheets.
$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 365
$ws1.Range("F4").Value  = 417
$ws1.Range("F5").Value  = 1132
$ws1.Range("F7").Value  = 44
$ws1.Range("F8").Value  = 944
$ws1.Range("F9").Value  = 1623
$ws1.Range("F10").Value = 6090
$ws1.Range("F11").Value = 112
$ws1.Range("F12").Value = 1750
$ws1.Range("F13").Value = 445
$ws1.Range("F14").Value = 5964
$ws1.Range("F18").Value = 96
$ws1.Range("F19").Value = 1655
$ws1.Range("F21").Value = 2
$ws1.Range("F22").Value = 37
$ws1.Range("F23").Value = 143
$ws1.Range("F24").Value = 1399
$ws1.Range("F25").Value = 724
$ws1.Range("F26").Value = 245

# Sheet: 演出
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value  = 165
$ws2.Range("F13").Value = 3

# Sheet: 本地生活
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 9501
$ws3.Range("F3").Value = 2233
$ws3.Range("F4").Value = 625
$ws3.Range("F5").Value = 189

# Sheet: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 9501
$ws4.Range("F3").Value  = 2233
$ws4.Range("F4").Value  = 625
$ws4.Range("F5").Value  = 365
$ws4.Range("F6").Value  = 417
$ws4.Range("F7").Value  = 1132
$ws4.Range("F9").Value  = 44
$ws4.Range("F12").Value = 944
$ws4.Range("F13").Value = 189
$ws4.Range("F14").Value = 1623
$ws4.Range("F15").Value = 6090
$ws4.Range("F16").Value = 112
$ws4.Range("F17").Value = 1750
$ws4.Range("F20").Value = 445
$ws4.Range("F23").Value = 5964
$ws4.Range("F27").Value = 96
$ws4.Range("F28").Value = 1654
$ws4.Range("F30").Value = 37
$ws4.Range("F31").Value = 143
$ws4.Range("F32").Value = 1399
$ws4.Range("F33").Value = 724
$ws4.Range("F34").Value = 3
$ws4.Range("F35").Value = 245
